# Fix a typo in the workbook's absolute path metadata:
# "ПоказателиЦУР" -> "Показатели ЦУР" (missing space between words)
$wb = $excel.ActiveWorkbook
$wb.AbsPath = "C:\Users\korozbaeva\Desktop\Показатели ЦУР для Платформы\Глобальные показатели ЦУР\"

$ws = $wb.ActiveSheet

# Update group-header labels on row 19 (age group) and row 29 (education group)
# to the "By ..." / "По ..." wording, and the Kyrgyz label to the longer form.
# Column order matters for how new shared strings get appended, so go
# column-by-column (A19, A29, B19, B29, C19, C29).
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"

$ws.Range("B19").Value = "По возрасту (в годах)"
$ws.Range("B29").Value = "По образованию"

$ws.Range("C19").Value = "By age (in years) "
$ws.Range("C29").Value = "By education"
